$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.792.83"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "1.870.63"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'301.04"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.5341"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("D8").Value = "'0.3743"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").Value = "'0.07186"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").Value = "'21.64"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "'0.08164"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "1.880.78"
$ws.Range("E13").Value = "  +25.72%  "
$ws.Range("D14").Value = "'92.89"
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'0.000008524"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "26.827.18"
$ws.Range("D21").Value = "'4.984"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "'2.315"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "'146.07"
$ws.Range("E25").Value = "  -2.80%  "
$ws.Range("D26").Value = "'1.734"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "'18.04"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "'113.98"
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("D29").Value = "'4.719"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").Value = "'4.631"
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").Value = "  -3.27%  "
$ws.Range("D33").Value = "'0.05027"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "'1.176"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("D35").Value = "'2.940"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").Value = "'0.6115"
$ws.Range("E36").Value = "  +5.79%  "
$ws.Range("D37").Value = "'2.694"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").Value = "'3.195"
$ws.Range("E38").Value = "  -4.74%  "
$ws.Range("D39").Value = "'0.01952"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("D40").Value = "'1.064"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.529"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("B42").Value = "Decentraland"
$ws.Range("C42").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D42").Value = "'0.5219"
$ws.Range("E42").Value = "  +6.11%  "
$ws.Range("D43").Value = "'8.771"
$ws.Range("E43").Value = "  -4.83%  "
$ws.Range("D44").Value = "'114.80"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "'1.646"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "'9.969"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").Value = "'37.62"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("D50").Value = "'0.06055"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'62.21"
$ws.Range("E51").Value = "  -3.41%  "
